# Normalize the "Recorded By" column (G) so that when a cell lists exactly
# two comma-separated recorders and one of them is "System", "System" is
# listed first, e.g. "dnasr281@gmail.com, System" -> "System, dnasr281@gmail.com".
# Cells that already start with "System", have only one recorder, or list
# more than two recorders are left untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)  # Column G = "Recorded By"
    $text = $cell.Text

    if ($text -and $text.Contains(", System")) {
        $parts = $text.Split(",")
        if ($parts.Count -eq 2) {
            $first = $parts[0].Trim()
            $second = $parts[1].Trim()
            if ($second -eq "System" -and $first -ne "System") {
                $cell.Value = "System, " + $first
            }
        }
    }
}
